# This script applies a cyclic permutation of the content of rows
# 4,5,6,7,8,9,11,12,13,14 on the "Artfynd" sheet (row 10 is left untouched).
# For each destination row, the values of columns A,B,D,E,F,G,H,M,Q,R,AO are
# replaced with the values that currently live in a different (source) row,
# per the mapping below. All reads are taken from the ORIGINAL values first,
# then all writes are applied, so that rows do not clobber each other while
# the permutation is being carried out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (source row's current values get copied
# into the destination row)
$mapping = @{
    4  = 9
    5  = 11
    6  = 4
    7  = 8
    8  = 13
    9  = 7
    11 = 6
    12 = 14
    13 = 12
    14 = 5
}

# Columns that participate in the row-content swap, plus their Excel
# column index (A=1 ... M=13 ... Q=17 ... R=18 ... AO=41).
$colIndex = @{
    "A"  = 1
    "B"  = 2
    "D"  = 4
    "E"  = 5
    "F"  = 6
    "G"  = 7
    "H"  = 8
    "M"  = 13
    "Q"  = 17
    "R"  = 18
    "AO" = 41
}

# 1. Snapshot the current (pre-edit) values of every involved column for
#    every row that is referenced (as either a source or destination).
$rowsInvolved = @{}
foreach ($dst in $mapping.Keys) { $rowsInvolved[$dst] = $true }
foreach ($src in $mapping.Values) { $rowsInvolved[$src] = $true }

$snapshot = @{}
foreach ($r in $rowsInvolved.Keys) {
    $rowData = @{}
    foreach ($colName in $colIndex.Keys) {
        $c = $colIndex[$colName]
        $rowData[$colName] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# 2. Apply the snapshot values from the source row into each destination
#    row. Column M is special: it is only populated for row 14 before the
#    edit (and must end up populated only for row 12 after the edit), so
#    clear it out first for all destination rows and only (re)write it
#    where the source actually had a value.
foreach ($dst in $mapping.Keys) {
    $src = $mapping[$dst]
    $data = $snapshot[$src]

    foreach ($colName in $colIndex.Keys) {
        $c = $colIndex[$colName]
        if ($colName -eq "M") {
            continue
        }
        $ws.Cells.Item($dst, $c).Value2 = $data[$colName]
    }

    # Handle column M explicitly: clear, then set only if the source had it.
    $mCol = $colIndex["M"]
    $ws.Cells.Item($dst, $mCol).Value2 = $null
    if ($null -ne $data["M"]) {
        $ws.Cells.Item($dst, $mCol).Value2 = $data["M"]
    }
}
